$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 131046938
$ws.Range("B7").Value = 79245
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = "Garnlav"
$ws.Range("G7").Value = "Alectoria sarmentosa"
$ws.Range("H7").Value = "(Ach.) Ach."
$ws.Range("M7").Value = ""
$ws.Range("Q7").Value = 395437
$ws.Range("R7").Value = 6804676
$ws.Range("Z7").Value = "11:35"
$ws.Range("AB7").Value = "11:35"
$ws.Range("A8").Value = 131046908
$ws.Range("Q8").Value = 395394
$ws.Range("R8").Value = 6804623
$ws.Range("Z8").Value = "11:02"
$ws.Range("AB8").Value = "11:02"
$ws.Range("A9").Value = 131046910
$ws.Range("Q9").Value = 395389
$ws.Range("R9").Value = 6804638
$ws.Range("Z9").Value = "11:06"
$ws.Range("AB9").Value = "11:06"
$ws.Range("A10").Value = 131046758
$ws.Range("B10").Value = 57881
$ws.Range("E10").Value = 100049
$ws.Range("F10").Value = "Spillkråka"
$ws.Range("G10").Value = "Dryocopus martius"
$ws.Range("H10").Value = "(Linnaeus, 1758)"
$ws.Range("M10").Value = "färska spår"
$ws.Range("Q10").Value = 395526
$ws.Range("R10").Value = 6804768
$ws.Range("Z10").Value = "12:08"
$ws.Range("AB10").Value = "12:08"
$ws.Range("A28").Value = 131046945
$ws.Range("B28").Value = 79245
$ws.Range("E28").Value = 6425
$ws.Range("F28").Value = "Garnlav"
$ws.Range("G28").Value = "Alectoria sarmentosa"
$ws.Range("H28").Value = "(Ach.) Ach."
$ws.Range("Q28").Value = 395448
$ws.Range("R28").Value = 6804579
$ws.Range("Z28").Value = "11:39"
$ws.Range("AB28").Value = "11:39"
$ws.Range("A29").Value = 131046991
$ws.Range("Q29").Value = 395479
$ws.Range("R29").Value = 6804588
$ws.Range("Z29").Value = "12:16"
$ws.Range("AB29").Value = "12:16"
$ws.Range("A30").Value = 131046937
$ws.Range("Q30").Value = 395433
$ws.Range("R30").Value = 6804692
$ws.Range("Z30").Value = "11:34"
$ws.Range("AB30").Value = "11:34"
$ws.Range("A31").Value = 131046964
$ws.Range("Q31").Value = 395381
$ws.Range("R31").Value = 6804728
$ws.Range("Z31").Value = "11:52"
$ws.Range("AB31").Value = "11:52"
$ws.Range("A32").Value = 131046914
$ws.Range("Q32").Value = 395382
$ws.Range("R32").Value = 6804689
$ws.Range("Z32").Value = "11:08"
$ws.Range("AB32").Value = "11:08"
$ws.Range("A33").Value = 131046801
$ws.Range("B33").Value = 81230
$ws.Range("E33").Value = 1049
$ws.Range("F33").Value = "Kortskaftad ärgspik"
$ws.Range("G33").Value = "Microcalicium ahlneri"
$ws.Range("H33").Value = "Tibell"
$ws.Range("Q33").Value = 395526
$ws.Range("R33").Value = 6804768
$ws.Range("Z33").Value = "12:08"
$ws.Range("AB33").Value = "12:08"
$ws.Range("A34").Value = 131046725
$ws.Range("B34").Value = 79277
$ws.Range("E34").Value = 185
$ws.Range("F34").Value = "Violettgrå tagellav"
$ws.Range("G34").Value = "Bryoria nadvornikiana"
$ws.Range("H34").Value = "(Gyeln.) Brodo & D.Hawksw."
$ws.Range("Q34").Value = 395392
$ws.Range("R34").Value = 6804595
$ws.Range("Z34").Value = "10:47"
$ws.Range("AB34").Value = "10:47"
$ws.Range("A35").Value = 131046987
$ws.Range("B35").Value = 79245
$ws.Range("E35").Value = 6425
$ws.Range("F35").Value = "Garnlav"
$ws.Range("G35").Value = "Alectoria sarmentosa"
$ws.Range("H35").Value = "(Ach.) Ach."
$ws.Range("Q35").Value = 395501
$ws.Range("R35").Value = 6804647
$ws.Range("Z35").Value = "12:14"
$ws.Range("AB35").Value = "12:14"
$ws.Range("A57").Value = 131046969
$ws.Range("B57").Value = 79245
$ws.Range("E57").Value = 6425
$ws.Range("F57").Value = "Garnlav"
$ws.Range("G57").Value = "Alectoria sarmentosa"
$ws.Range("H57").Value = "(Ach.) Ach."
$ws.Range("Q57").Value = 395464
$ws.Range("R57").Value = 6804807
$ws.Range("Z57").Value = "11:57"
$ws.Range("AB57").Value = "11:57"
$ws.Range("A58").Value = 131046926
$ws.Range("Q58").Value = 395395
$ws.Range("R58").Value = 6804786
$ws.Range("Z58").Value = "11:26"
$ws.Range("AB58").Value = "11:26"
$ws.Range("A59").Value = 131046925
$ws.Range("Q59").Value = 395380
$ws.Range("R59").Value = 6804774
$ws.Range("Z59").Value = "11:25"
$ws.Range("AB59").Value = "11:25"
$ws.Range("A60").Value = 131046722
$ws.Range("B60").Value = 79277
$ws.Range("E60").Value = 185
$ws.Range("F60").Value = "Violettgrå tagellav"
$ws.Range("G60").Value = "Bryoria nadvornikiana"
$ws.Range("H60").Value = "(Gyeln.) Brodo & D.Hawksw."
$ws.Range("Q60").Value = 395391
$ws.Range("R60").Value = 6804603
$ws.Range("Z60").Value = "10:52"
$ws.Range("AB60").Value = "10:52"
$ws.Range("A61").Value = 131046904
$ws.Range("Q61").Value = 395385
$ws.Range("R61").Value = 6804578
$ws.Range("Z61").Value = "10:43"
$ws.Range("AB61").Value = "10:43"
$ws.Range("A62").Value = 131046939
$ws.Range("Q62").Value = 395446
$ws.Range("R62").Value = 6804659
$ws.Range("Z62").Value = "11:35"
$ws.Range("AB62").Value = "11:35"
$ws.Range("A63").Value = 131046983
$ws.Range("Q63").Value = 395515
$ws.Range("R63").Value = 6804694
$ws.Range("Z63").Value = "12:12"
$ws.Range("AB63").Value = "12:12"
$ws.Range("A68").Value = 131047026
$ws.Range("B68").Value = 78911
$ws.Range("E68").Value = 353
$ws.Range("F68").Value = "Dvärgbägarlav"
$ws.Range("G68").Value = "Cladonia parasitica"
$ws.Range("H68").Value = "(Hoffm.) Hoffm."
$ws.Range("Q68").Value = 395386
$ws.Range("R68").Value = 6804628
$ws.Range("Z68").Value = "11:04"
$ws.Range("AB68").Value = "11:04"
$ws.Range("A70").Value = 131046917
$ws.Range("B70").Value = 79245
$ws.Range("E70").Value = 6425
$ws.Range("F70").Value = "Garnlav"
$ws.Range("G70").Value = "Alectoria sarmentosa"
$ws.Range("H70").Value = "(Ach.) Ach."
$ws.Range("Q70").Value = 395365
$ws.Range("R70").Value = 6804704
$ws.Range("Z70").Value = "11:11"
$ws.Range("AB70").Value = "11:11"
$ws.Range("A71").Value = 131046921
$ws.Range("Q71").Value = 395357
$ws.Range("R71").Value = 6804752
$ws.Range("Z71").Value = "11:19"
$ws.Range("AB71").Value = "11:19"
$ws.Range("A73").Value = 131047012
$ws.Range("B73").Value = 57884
$ws.Range("E73").Value = 100109
$ws.Range("F73").Value = "Tretåig hackspett"
$ws.Range("G73").Value = "Picoides tridactylus"
$ws.Range("H73").Value = "(Linnaeus, 1758)"
$ws.Range("M73").Value = "färska spår"
$ws.Range("R73").Value = 6804659
$ws.Range("Z73").Value = "11:35"
$ws.Range("AB73").Value = "11:35"
$ws.Range("AC73").Value = "Troliga spår efter tretåig hackspett (barkfälkning)"
$ws.Range("AE73").Value = $true
$ws.Range("A74").Value = 131046930
$ws.Range("Q74").Value = 395446
$ws.Range("R74").Value = 6804802
$ws.Range("Z74").Value = "11:29"
$ws.Range("AB74").Value = "11:29"
$ws.Range("A75").Value = 131046916
$ws.Range("Q75").Value = 395367
$ws.Range("R75").Value = 6804698
$ws.Range("Z75").Value = "11:09"
$ws.Range("AB75").Value = "11:09"
$ws.Range("A76").Value = 131046933
$ws.Range("B76").Value = 79245
$ws.Range("E76").Value = 6425
$ws.Range("F76").Value = "Garnlav"
$ws.Range("G76").Value = "Alectoria sarmentosa"
$ws.Range("H76").Value = "(Ach.) Ach."
$ws.Range("M76").Value = ""
$ws.Range("Q76").Value = 395458
$ws.Range("R76").Value = 6804762
$ws.Range("Z76").Value = "11:32"
$ws.Range("AB76").Value = "11:32"
$ws.Range("AC76").Value = ""
$ws.Range("AE76").Value = $false
$ws.Range("A91").Value = 131046804
$ws.Range("B91").Value = 79002
$ws.Range("E91").Value = 6446
$ws.Range("F91").Value = "Kolflarnlav"
$ws.Range("G91").Value = "Carbonicola anthracophila"
$ws.Range("H91").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("M91").Value = ""
$ws.Range("Q91").Value = 395367
$ws.Range("R91").Value = 6804754
$ws.Range("Z91").Value = "11:21"
$ws.Range("AB91").Value = "11:21"
$ws.Range("AC91").Value = ""
$ws.Range("A92").Value = 131046718
$ws.Range("B92").Value = 83225
$ws.Range("E92").Value = 6440
$ws.Range("F92").Value = "Vitgrynig nållav"
$ws.Range("G92").Value = "Chaenotheca subroscida"
$ws.Range("H92").Value = "(Eitner) Zahlbr."
$ws.Range("Q92").Value = 395471
$ws.Range("R92").Value = 6804595
$ws.Range("Z92").Value = "12:16"
$ws.Range("AB92").Value = "12:16"
$ws.Range("A93").Value = 131046984
$ws.Range("B93").Value = 79245
$ws.Range("E93").Value = 6425
$ws.Range("F93").Value = "Garnlav"
$ws.Range("G93").Value = "Alectoria sarmentosa"
$ws.Range("H93").Value = "(Ach.) Ach."
$ws.Range("Q93").Value = 395531
$ws.Range("R93").Value = 6804688
$ws.Range("Z93").Value = "12:12"
$ws.Range("AB93").Value = "12:12"
$ws.Range("A94").Value = 131046791
$ws.Range("B94").Value = 57884
$ws.Range("E94").Value = 100109
$ws.Range("F94").Value = "Tretåig hackspett"
$ws.Range("G94").Value = "Picoides tridactylus"
$ws.Range("H94").Value = "(Linnaeus, 1758)"
$ws.Range("M94").Value = "färska spår"
$ws.Range("Q94").Value = 395362
$ws.Range("R94").Value = 6804701
$ws.Range("Z94").Value = "11:10"
$ws.Range("AB94").Value = "11:10"
$ws.Range("AC94").Value = "Färska ringhack (gran)"
